$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values entered by the user into C2 and C3
$ws.Range("C2").Value = 4
$ws.Range("C3").Value = 2

# C4 used to hold a SUM formula; it is now overwritten with a plain value
$ws.Range("C4").Value = 7

# A new running-total formula is entered in C5, recalculated after the
# values above were entered
$ws.Range("C5").Formula = "=SUM(C2:C4)"

# Move the active selection
$ws.Range("N16").Select()
